# edit.ps1 - apply the "Quantum Entanglement" -> "Biology" rewrite described
# by the target diff, using Word COM-interop calls against $word.ActiveDocument.
#
# Strategy: use Find.Execute purely to *locate* text (MatchCase=$true, no
# auto-replace), then assign the new text directly to the located Range's
# .Text property. Doing the substitution this way (rather than passing a
# ReplaceWith string straight into Find.Execute) avoids the engine's
# "AutoCorrect while typing" smart-quote mangling, so straight apostrophes
# stay straight.

function Replace-ExactText {
    param(
        $ScopeRange,
        [string]$Old,
        [string]$New
    )
    $rng = $ScopeRange.Duplicate
    $found = $rng.Find.Execute($Old, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Replace-ExactText: text not found -> " + $Old)
    }
    $rng.Text = $New
    return $rng
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------
Replace-ExactText $d.Paragraphs(1).Range `
    "Quantum Entanglement: A Mysterious Dance" `
    "The Captivating World of Biology: Unveiling the Secrets of Life" | Out-Null

# ---------------------------------------------------------------------
# 2. Author line: "Dr" + "." + " Anna Smith" (3 runs) -> "Olivia Rose"
# ---------------------------------------------------------------------
Replace-ExactText $d.Paragraphs(2).Range "Dr. Anna Smith" "Olivia Rose" | Out-Null

# ---------------------------------------------------------------------
# 3. Email line: "anna" / "smith@quantuminstitute" / "org" -> new values
#    (the "." separator runs in between are left alone)
# ---------------------------------------------------------------------
Replace-ExactText $d.Paragraphs(3).Range "anna" "olivia" | Out-Null
Replace-ExactText $d.Paragraphs(3).Range "smith@quantuminstitute" "rose@validschool" | Out-Null
Replace-ExactText $d.Paragraphs(3).Range "org" "edu" | Out-Null

# ---------------------------------------------------------------------
# 4. Body paragraph (quantum-mechanics blurb -> biology blurb)
# ---------------------------------------------------------------------
$bodyPairs = @(
    @("The realm of quantum mechanics is a fascinating and enigmatic frontier of science, pushing the boundaries of our understanding of the universe",
      "Embark on an enthralling journey into the realm of biology, where we unravel the intricate tapestry of life"),
    @(" Among its many intriguing phenomena, quantum entanglement stands out as one of the most mysterious and counterintuitive",
      " Delve into the fascinating microscopic world of cells, the building blocks of all living organisms, and discover their remarkable capabilities"),
    @(' Einstein famously referred to it as "spooky action at a distance," and its implications continue to perplex and enthrall scientists and philosophers alike',
      ' Explore the intricate processes that govern genetics, the blueprint of life, and witness the mesmerizing dance of molecules, the fundamental components of all matter'),
    @("This mystical dance between particles defies classical notions of locality and causality, allowing them to share information instantaneously, regardless of the distance separating them",
      "Journey through the awe-inspiring diversity of life on Earth, from the towering giants of the rainforest to the microscopic organisms that thrive in extreme environments"),
    @(" It's as if they are connected by an invisible thread, responding to each other's actions in perfect synchrony, even across vast cosmic distances",
      " Witness the intricate adaptations that enable organisms to survive and thrive in a myriad of habitats, showcasing the resilience and adaptability of life"),
    @(" The implications of this phenomenon are profound, challenging our understanding of reality and opening up new possibilities in the realm of communication, computing, and cryptography",
      " Delve into the delicate balance of ecosystems, where organisms interact in complex webs of interdependence, understanding the profound impact of human activities on the delicate balance of nature"),
    @("Quantum entanglement has been experimentally verified numerous times, leaving no doubt about its existence",
      "Unravel the mysteries of the human body, a marvel of engineering, and explore the intricate workings of its systems"),
    @(" However, the underlying mechanism responsible for this strange phenomenon remains a subject of intense debate and research",
      " Discover the intricate network of organs, tissues, and cells that orchestrate a symphony of life, maintaining homeostasis and enabling us to interact with the world around us"),
    @(" Some physicists posit the existence of hidden variables that govern the behavior of entangled particles, while others propose that the particles themselves possess a form of non-local consciousness, allowing them to communicate instantaneously",
      " Investigate the fascinating processes of growth, development, and reproduction, marveling at the intricacies of life's creation")
)

# body text lives in paragraph 5
foreach ($pair in $bodyPairs) {
    Replace-ExactText $d.Paragraphs(5).Range $pair[0] $pair[1] | Out-Null
}

# ---------------------------------------------------------------------
# 5. Summary paragraph: 3 sentences rewritten, trailing sentences removed
# ---------------------------------------------------------------------
# Summary is the last paragraph of the body text.
$summaryIndex = $d.Paragraphs.Count

Replace-ExactText $d.Paragraphs($summaryIndex).Range `
    "Quantum entanglement is a mysterious phenomenon in which particles become interconnected, sharing information instantaneously regardless of the distance separating them" `
    "Biology, the study of life, is an enthralling field that unveils the intricate tapestry of life on Earth" | Out-Null

Replace-ExactText $d.Paragraphs($summaryIndex).Range `
    " It defies classical notions of locality and causality, challenging our understanding of reality" `
    " From the microscopic realm of cells to the awe-inspiring diversity of organisms, biology offers a profound understanding of the fundamental processes that govern life" | Out-Null

# Drop everything from the middle of the "Despite ..." sentence through to
# the end of the "... nature of the universe" sentence, leaving the final
# "." run untouched and in place.
$summaryTail = ". The implications of quantum entanglement are profound, with potential applications in communication, computing, and cryptography." + `
    " As we delve deeper into the enigmatic world of quantum mechanics, unraveling the secrets of quantum entanglement may provide unprecedented insights into the fundamental nature of the universe"
Replace-ExactText $d.Paragraphs($summaryIndex).Range $summaryTail "" | Out-Null

Replace-ExactText $d.Paragraphs($summaryIndex).Range `
    " Despite numerous experimental verifications, the underlying mechanism responsible for this strange phenomenon remains a subject of intense debate" `
    " Exploring genetics, adaptations, ecosystems, and the human body, we gain invaluable insights into the wonders of the natural world, fostering a deeper appreciation for the delicate balance of life" | Out-Null

# ---------------------------------------------------------------------
# 6. Trailing empty paragraph added at the end of the document body
# ---------------------------------------------------------------------
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter() | Out-Null
